$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 30 holds the "dividendos" entry for the shares screen.
# Columns A/B/C mirror the formatting of the preceding data row (A29:C29)
# exactly (border + fill style), so copy that row down first, then fill in
# the values (in the same order they were typed when the row was authored).
$ws.Range("A29:C29").Copy($ws.Range("A30:C30"))

$ws.Range("C30").Value = "SP_DIVIDENDOS_SELECT"
$ws.Range("A30").Value = "dividendos"
$ws.Range("B30").Value = "DividendosController"

# G30/H30 use the shaded "Página"/"metodo" style (fill, no border) used
# elsewhere in the table; copy it from an existing shaded cell then strip
# the border so only the fill remains.
$ws.Range("G23").Copy($ws.Range("G30"))
$ws.Range("G23").Copy($ws.Range("H30"))
$ws.Range("G30:H30").Borders.LineStyle = -4142

$ws.Range("H30").Value = "ObtenerDividendos"
$ws.Range("G30").Value = "shares-page"

# D30 is a plain, unstyled cell.
$ws.Range("D30").Value = "emisor"

# Move the active selection the way the saved workbook records it.
$ws.Range("D31").Select()
